$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hardware Hours (column B, e.g. "0a444b4") together with the Author
# (column A) are now also populated for the "Major Componnets BOM" task
# rows (26-29), matching the pattern already used for every other task
# row in the table. Row 26 already had its category (column C) filled
# in; rows 27-29 need that too.
foreach ($r in 26..29) {
    $ws.Range("A25").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("A$r").Value2 = $ws.Range("A25").Value2

    $ws.Range("B25").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("B$r").Value2 = $ws.Range("B25").Value2

    if ($r -ne 26) {
        $ws.Range("C26").Copy()
        $ws.Range("C$r").PasteSpecial(-4122)   # xlPasteFormats
        $ws.Range("C$r").Value2 = $ws.Range("C26").Value2
    }
}

# A new trailing row is added, carrying forward the column-B formatting
# only (no value), matching the end of the enumerated task list.
$ws.Range("B25").Copy()
$ws.Range("B30").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# Leave the selection where the author ended up after the edit.
$ws.Range("B27").Select()
